$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the "Price" column that are being updated. Several of the new
# prices (e.g. "608.63") look like plain numbers, so Excel would silently
# convert them to numeric values on assignment. Force each cell to Text
# format first so the values are stored as strings, matching the source data.
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated price values
$ws.Range("D2").Value = "69.315.39"
$ws.Range("D3").Value = "3.492.40"
$ws.Range("D5").Value = "608.63"
$ws.Range("D6").Value = "185.94"
$ws.Range("D9").Value = "0.210"
$ws.Range("D10").Value = "0.651"
$ws.Range("D11").Value = "53.34"
$ws.Range("D12").Value = "0.0000306"
$ws.Range("D14").Value = "4.049.27"
$ws.Range("D15").Value = "608.19"
$ws.Range("D16").Value = "18.99"
$ws.Range("D17").Value = "12.69"
$ws.Range("D18").Value = "69.377.37"
$ws.Range("D19").Value = "3.497.28"
$ws.Range("D22").Value = "17.39"
$ws.Range("D23").Value = "103.97"
$ws.Range("D24").Value = "4.64"
$ws.Range("D25").Value = "5.04"
$ws.Range("D26").Value = "3.04"
$ws.Range("D27").Value = "10.93"
$ws.Range("D28").Value = "9.92"
$ws.Range("D29").Value = "33.64"
$ws.Range("D30").Value = "7.02"
$ws.Range("D31").Value = "12.51"
$ws.Range("D33").Value = "63.31"
$ws.Range("D34").Value = "3.78"
$ws.Range("D35").Value = "3.15"
$ws.Range("D36").Value = "0.999"
$ws.Range("D37").Value = "524.95"
$ws.Range("D38").Value = "0.396"
$ws.Range("D39").Value = "3.59"
$ws.Range("D40").Value = "3.555.58"
$ws.Range("D41").Value = "36.67"
$ws.Range("D42").Value = "0.0₃0767"
$ws.Range("D43").Value = "0.139"
$ws.Range("D44").Value = "0.0461"
$ws.Range("D45").Value = "2.99"
$ws.Range("D47").Value = "3.31"
$ws.Range("D48").Value = "8.87"
$ws.Range("D50").Value = "131.58"

# The values are now stored as text; restore the default (General) cell
# formatting so no stray number format is left behind on these cells.
foreach ($addr in $priceCells) {
    $ws.Range($addr).ClearFormats()
}

# Updated coin name / link / 1h-volume values. These are never ambiguous with
# numbers (names, URLs, or "  +/-X.XX%  " strings), so no format fixing is needed.
$ws.Range("E2").Value = "  -2.28%  "
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  +4.38%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -2.88%  "
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("E12").Value = "  -3.84%  "
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("E15").Value = "  +8.01%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("E19").Value = "  -1.87%  "
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("E21").Value = "  -2.27%  "
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("E23").Value = "  +9.94%  "
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  +3.00%  "
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("E28").Value = "  +8.08%  "
$ws.Range("E29").Value = "  +3.45%  "
$ws.Range("E30").Value = "  -3.95%  "
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E34").Value = "  +14.21%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E35").Value = "  -7.75%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  -5.15%  "
$ws.Range("E38").Value = "  -5.65%  "
$ws.Range("E39").Value = "  +4.61%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("E42").Value = "  -4.40%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("E44").Value = "  +2.55%  "
$ws.Range("E45").Value = "  +1.68%  "
$ws.Range("E46").Value = "  +4.18%  "
$ws.Range("E47").Value = "  -5.15%  "
$ws.Range("E48").Value = "  -5.57%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("E51").Value = "  -9.07%  "
